$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DEMANDS")
$ws.Range("I22").Value = 1.23
